$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 116 ("Politiezone onbekend (Vlaanderen)")
# so it shifts down to row 117, and a new row 116 is created for the
# "Politiezone onbekend (Brussel)" entry (previously row 115).
$ws.Rows.Item(116).Insert()

# Row 115 now becomes "Politiezone onbekend" (no suffix), code ipz993.
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = "ipz993"
$ws.Cells.Item(115, 3).Value = "Politiezone onbekend"
$ws.Cells.Item(115, 4).Value = "Politiezone onbekend"

# New row 116 gets the old "Politiezone onbekend (Brussel)" entry, code ipz992.
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "ipz992"
$ws.Cells.Item(116, 3).Value = "Politiezone onbekend (Brussel)"
$ws.Cells.Item(116, 4).Value = "Politiezone onbekend (Brussel)"

# Row 117 (shifted from old row 116) keeps "Politiezone onbekend (Vlaanderen)", code ipz991,
# just update the volgnr to 116 to match the new sequence.
$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = "ipz991"
$ws.Cells.Item(117, 3).Value = "Politiezone onbekend (Vlaanderen)"
$ws.Cells.Item(117, 4).Value = "Politiezone onbekend (Vlaanderen)"
